# Updates the cryptos Price (D) and Volume(1h) (E) columns with the
# latest scraped values. Plain decimal-looking Price values are forced
# to stay text (matching the source sheet's inline-string cells) by
# temporarily applying a text NumberFormat before writing them, then
# restoring the "Normal" style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "47.909.78"
$ws.Range("E2").Value = "  -0.53%  "

# Row 3
$ws.Range("D3").Value = "2.486.46"
$ws.Range("E3").Value = "  -1.13%  "

# Row 4
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.16%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.37%  "

# Row 7
$ws.Range("E7").Value = "  -2.35%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.52%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.29%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.24"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.17%  "

# Row 12
$ws.Range("E12").Value = "  -2.57%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.125"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.55%  "

# Row 14
$ws.Range("E14").Value = "  -2.64%  "

# Row 15
$ws.Range("D15").Value = "2.874.46"
$ws.Range("E15").Value = "  -1.25%  "

# Row 16
$ws.Range("D16").Value = "2.463.31"
$ws.Range("E16").Value = "  -2.18%  "

# Row 17
$ws.Range("E17").Value = "  -2.50%  "

# Row 18
$ws.Range("D18").Value = "47.776.25"
$ws.Range("E18").Value = "  -0.49%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.35%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.79%  "

# Row 21
$ws.Range("E21").Value = "  -1.22%  "

# Row 22
$ws.Range("E22").Value = "  -2.21%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "281.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.94%  "

# Row 24
$ws.Range("E24").Value = "  -1.85%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.28%  "

# Row 26
$ws.Range("E26").Value = "  -0.10%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.44%  "

# Row 28
$ws.Range("E28").Value = "  -6.78%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.63%  "

# Row 30
$ws.Range("E30").Value = "  -3.62%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.24%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.06%  "

# Row 33
$ws.Range("E33").Value = "  -0.15%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.92%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0767"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.13%  "

# Row 37
$ws.Range("E37").Value = "  -2.01%  "

# Row 38
$ws.Range("E38").Value = "  -4.31%  "

# Row 39
$ws.Range("E39").Value = "  -3.62%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.111"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.98%  "

# Row 41
$ws.Range("E41").Value = "  -1.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "119.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.79%  "

# Row 43
$ws.Range("E43").Value = "  -2.42%  "

# Row 44
$ws.Range("E44").Value = "  -0.81%  "

# Row 45
$ws.Range("D45").Value = "1.983.46"
$ws.Range("E45").Value = "  -2.35%  "

# Row 46
$ws.Range("E46").Value = "  -1.08%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.27%  "

# Row 48
$ws.Range("E48").Value = "  +2.32%  "

# Row 49
$ws.Range("E49").Value = "  -1.63%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.13%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.42%  "
